$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $val) {
    $q = [string][char]39
    $sheet.Range($addr).Value = $q + $val
    $sheet.Range($addr).Style = "Normal"
}

Set-CellText $ws 'D2' '30.197.62'
Set-CellText $ws 'E2' '  -0.39%  '
Set-CellText $ws 'D3' '1.862.69'
Set-CellText $ws 'E3' '  -0.36%  '
Set-CellText $ws 'D4' '0.9990'
Set-CellText $ws 'E4' '  -0.25%  '
Set-CellText $ws 'D5' '243.35'
Set-CellText $ws 'E5' '  +3.50%  '
Set-CellText $ws 'D6' '0.9994'
Set-CellText $ws 'E6' '  -0.17%  '
Set-CellText $ws 'D7' '0.4726'
Set-CellText $ws 'E7' '  +0.68%  '
Set-CellText $ws 'D8' '42.80'
Set-CellText $ws 'E8' '  -0.32%  '
Set-CellText $ws 'D9' '0.2857'
Set-CellText $ws 'E9' '  -0.22%  '
Set-CellText $ws 'D10' '0.06476'
Set-CellText $ws 'E10' '  -1.75%  '
Set-CellText $ws 'D12' '0.07682'
Set-CellText $ws 'E12' '  -3.81%  '
Set-CellText $ws 'D13' '1.861.93'
Set-CellText $ws 'E13' '  -0.50%  '
Set-CellText $ws 'D14' '94.31'
Set-CellText $ws 'E14' '  -2.60%  '
Set-CellText $ws 'B15' 'Polkadot'
Set-CellText $ws 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-CellText $ws 'D15' '5.076'
Set-CellText $ws 'E15' '  -0.69%  '
Set-CellText $ws 'B16' 'Polygon'
Set-CellText $ws 'C16' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-CellText $ws 'D16' '0.6821'
Set-CellText $ws 'E16' '  -1.03%  '
Set-CellText $ws 'D17' '270.84'
Set-CellText $ws 'E17' '  +0.73%  '
Set-CellText $ws 'D18' '30.189.31'
Set-CellText $ws 'D19' '13.36'
Set-CellText $ws 'E19' '  -5.60%  '
Set-CellText $ws 'D20' '0.000007555'
Set-CellText $ws 'E20' '  -2.67%  '
Set-CellText $ws 'E21' '  -0.16%  '
Set-CellText $ws 'D22' '2.118.68'
Set-CellText $ws 'E22' '  +0.13%  '
Set-CellText $ws 'D23' '0.9986'
Set-CellText $ws 'E23' '  -0.15%  '
Set-CellText $ws 'D24' '5.181'
Set-CellText $ws 'E24' '  -1.47%  '
Set-CellText $ws 'D25' '6.105'
Set-CellText $ws 'E25' '  -1.76%  '
Set-CellText $ws 'D26' '9.342'
Set-CellText $ws 'E26' '  -0.53%  '
Set-CellText $ws 'D27' '165.81'
Set-CellText $ws 'E27' '  -0.96%  '
Set-CellText $ws 'D28' '18.78'
Set-CellText $ws 'E28' '  -0.54%  '
Set-CellText $ws 'D29' '1.887'
Set-CellText $ws 'E29' '  -3.15%  '
Set-CellText $ws 'D30' '1.374'
Set-CellText $ws 'E30' '  +0.73%  '
Set-CellText $ws 'D31' '0.09854'
Set-CellText $ws 'E31' '  -0.14%  '
Set-CellText $ws 'D32' '1.508'
Set-CellText $ws 'E32' '  +3.44%  '
Set-CellText $ws 'D33' '4.235'
Set-CellText $ws 'E33' '  -2.63%  '
Set-CellText $ws 'D34' '3.985'
Set-CellText $ws 'E34' '  -1.79%  '
Set-CellText $ws 'D35' '0.04710'
Set-CellText $ws 'E35' '  -0.01%  '
Set-CellText $ws 'E36' '  -2.20%  '
Set-CellText $ws 'D37' '0.6873'
Set-CellText $ws 'E37' '  -2.19%  '
Set-CellText $ws 'D38' '2.706'
Set-CellText $ws 'E38' '  -1.23%  '
Set-CellText $ws 'D39' '0.01846'
Set-CellText $ws 'E39' '  -1.77%  '
Set-CellText $ws 'D40' '2.726'
Set-CellText $ws 'E40' '  -3.35%  '
Set-CellText $ws 'D41' '6.397'
Set-CellText $ws 'E41' '  +2.36%  '
Set-CellText $ws 'D42' '70.39'
Set-CellText $ws 'E42' '  -2.33%  '
Set-CellText $ws 'D43' '0.9988'
Set-CellText $ws 'E43' '  -0.14%  '
Set-CellText $ws 'D44' '0.8357'
Set-CellText $ws 'E44' '  -0.79%  '
Set-CellText $ws 'D45' '1.890'
Set-CellText $ws 'E45' '  -3.50%  '
Set-CellText $ws 'D46' '102.07'
Set-CellText $ws 'E46' '  -1.02%  '
Set-CellText $ws 'D47' '0.4069'
Set-CellText $ws 'E47' '  -2.55%  '
Set-CellText $ws 'D48' '9.201'
Set-CellText $ws 'E48' '  +0.36%  '
Set-CellText $ws 'D49' '929.27'
Set-CellText $ws 'E49' '  +1.16%  '
Set-CellText $ws 'D50' '6.951'
Set-CellText $ws 'E50' '  -1.94%  '
Set-CellText $ws 'D51' '34.47'
Set-CellText $ws 'E51' '  -0.12%  '

Write-Output "Applied 99 cell updates"
